$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Batch")
$ws.Activate()

# --- New header cells (row 1) ---
# D1: reuse the "Calibri" style used by B1/C1
$ws.Range("D1").Font.Name = "Calibri"
$ws.Range("D1").Value = "BatchNoOfClassesmissingadditonalfield"

# E1: same Calibri style
$ws.Range("E1").Font.Name = "Calibri"
$ws.Range("E1").Value = "BatchStatusmissingadditonalfield"

# F1: same Aptos style used by A1
$ws.Range("F1").Font.Name = "Aptos"
$ws.Range("F1").Value = "BatchDescriptionmissingadditionalfield"

# G1: empty cell, but still styled like the Calibri cells
$ws.Range("G1").Font.Name = "Calibri"

# --- New data cells (row 2) ---
$ws.Range("D2").Value = 4
$ws.Range("E2").Value = "Active"
$ws.Range("F2").Value = "null"

# --- Column widths for D:F ---
$ws.Columns.Item(4).ColumnWidth = 37.666666666666664
$ws.Columns.Item(5).ColumnWidth = 29.833333333333332
$ws.Columns.Item(6).ColumnWidth = 34.666666666666664

# --- Selection moves to E1 ---
$ws.Range("E1").Select()
